# Update Group Test Case: change the "Environment" column (P2:P11) from
# "OS: Windows 8.1 / Browser: Chrome 41" to "OS: Windows 7 / IDE: Eclipse Luna"
# and update the sheet's view/selection to reflect the edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEnv = "OS: Windows 7" + [char]10 + "IDE: Eclipse Luna"

foreach ($r in 2..11) {
    $ws.Cells.Item($r, 16).Value = $newEnv
}

# Reflect the selection / scroll position used while making this edit.
$ws.Range("P2:P11").Select()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 10
